$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Data")

# Helper dates used throughout the "BeginDateTime" / "EndDateTime" columns.
$beginDate = Get-Date -Year 2000 -Month 1 -Day 1 -Hour 0 -Minute 0 -Second 0
$endDate   = Get-Date -Year 9998 -Month 12 -Day 31 -Hour 0 -Minute 0 -Second 0

# ---------------------------------------------------------------------------
# Existing rows: extend the EndDateTime (column H) from 2099-01-01 to the
# "no expiry" sentinel 9998-12-31 for every already-present setting row.
# ---------------------------------------------------------------------------
$ws.Range("H2").Value = $endDate
$ws.Range("H3").Value = $endDate
$ws.Range("H4").Value = $endDate

# Re-affirm the (already blank) DisplayName / Comment cells on the untouched
# rows so they stay genuinely empty.
$ws.Range("D2").Value = ""
$ws.Range("F2").Value = ""
$ws.Range("D3").Value = ""
$ws.Range("F3").Value = ""
$ws.Range("D4").Value = ""
$ws.Range("F4").Value = ""

# Row 3 (Id = 2): rename setting from CurrentLpuName to OrgName and update its
# value/description to reflect the full legal organisation name.
$ws.Range("B3").Value = "OrgName"
$ws.Range("C3").Value = "Государственное бюджетное образовательное учреждение высшего профессионального образования «Российский национальный исследовательский медицинский университет имени Н.И.Пирогова» Министерства здравоохранения Российской Федерации"
$ws.Range("E3").Value = "Полное название ЛПУ"

# Row 4 (NotificationServiceAddress): its numeric Id moves from 3 to 5 to make
# room for the newly inserted settings below.
$ws.Range("A4").Value = "'5"

# ---------------------------------------------------------------------------
# New configuration rows (5-14).
# ---------------------------------------------------------------------------
$ws.Range("A5").Value = "'7"
$ws.Range("B5").Value = "OrgOKPO"
$ws.Range("C5").Value = "'11223444"
$ws.Range("D5").Value = ""
$ws.Range("E5").Value = "ОКПО ЛПУ"
$ws.Range("F5").Value = ""
$ws.Range("G5").Value = $beginDate
$ws.Range("H5").Value = $endDate

$ws.Range("A6").Value = "'8"
$ws.Range("B6").Value = "OrgShortName"
$ws.Range("C6").Value = "ГБОУ ВПО РНИМУ им. Н.И. Пирогова МЗ РФ"
$ws.Range("D6").Value = ""
$ws.Range("E6").Value = "Сокращенное название ЛПУ"
$ws.Range("F6").Value = ""
$ws.Range("G6").Value = $beginDate
$ws.Range("H6").Value = $endDate

$ws.Range("A7").Value = "'10"
$ws.Range("B7").Value = "NIKIName"
$ws.Range("C7").Value = "Научно-исследовательский клинический институт педиатрии им. академика Ю.Е. Вельтищева ГБОУ ВПО РНИМУ им. Н.И. Пирогова МЗ РФ"
$ws.Range("D7").Value = "Научно-исследовательском клиническом институте педиатрии им. академика Ю.Е. Вельтищева ГБОУ ВПО РНИМУ им. Н.И. Пирогова МЗ РФ"
$ws.Range("E7").Value = "Название обособленного структурного подразделения НИКИ"
$ws.Range("F7").Value = ""
$ws.Range("G7").Value = $beginDate
$ws.Range("H7").Value = $endDate

$ws.Range("A8").Value = "'11"
$ws.Range("B8").Value = "DirectorFullName"
$ws.Range("C8").Value = "Школьникова Мария Александровна"
$ws.Range("D8").Value = "Школьниковой Марии Александровны"
$ws.Range("E8").Value = "Директор НИКИ"
$ws.Range("F8").Value = ""
$ws.Range("G8").Value = $beginDate
$ws.Range("H8").Value = $endDate

$ws.Range("A9").Value = "'12"
$ws.Range("B9").Value = "PayContractLicense"
$ws.Range("C9").Value = "17 марта 2014 года № 52"
$ws.Range("D9").Value = ""
$ws.Range("E9").Value = "Доверенность на оказание платных услуг"
$ws.Range("F9").Value = ""
$ws.Range("G9").Value = $beginDate
$ws.Range("H9").Value = $endDate

$ws.Range("A10").Value = "'13"
$ws.Range("B10").Value = "NIKIAddress"
$ws.Range("C10").Value = "125412, г. Москва, ул. Талдомская, д. 2"
$ws.Range("D10").Value = ""
$ws.Range("E10").Value = "Адрес НИКИ"
$ws.Range("F10").Value = ""
$ws.Range("G10").Value = $beginDate
$ws.Range("H10").Value = $endDate

$ws.Range("A11").Value = "'14"
$ws.Range("B11").Value = "OrgAddress"
$ws.Range("C11").Value = "117997, г. Москва, ул. Островитянова, д. 1"
$ws.Range("D11").Value = ""
$ws.Range("E11").Value = "Юридический адрес"
$ws.Range("F11").Value = ""
$ws.Range("G11").Value = $beginDate
$ws.Range("H11").Value = $endDate

$ws.Range("A12").Value = "'15"
$ws.Range("B12").Value = "DirectorShortName"
$ws.Range("C12").Value = "Школьникова М.А."
$ws.Range("D12").Value = "Школьникова М.А."
$ws.Range("E12").Value = "Директор НИКИ"
$ws.Range("F12").Value = ""
$ws.Range("G12").Value = $beginDate
$ws.Range("H12").Value = $endDate

$ws.Range("A13").Value = "'16"
$ws.Range("B13").Value = "NIKIShortName"
$ws.Range("C13").Value = "Научно-исследовательский клинический институт педиатрии"
$ws.Range("D13").Value = "Научно-исследовательском клиническом институте педиатрии"
$ws.Range("E13").Value = "Сокращенное название обособленного структурного подразделения НИКИ"
$ws.Range("F13").Value = ""
$ws.Range("G13").Value = $beginDate
$ws.Range("H13").Value = $endDate

$ws.Range("A14").Value = "'18"
$ws.Range("B14").Value = "ChildAge"
$ws.Range("C14").Value = "'14"
$ws.Range("D14").Value = ""
$ws.Range("E14").Value = "Возраст, старше которого человек считается взрослым"
$ws.Range("F14").Value = ""
$ws.Range("G14").Value = $beginDate
$ws.Range("H14").Value = $endDate
